# Auto-generated edit script applying numeric updates to the Leviathan_Profits workbook
# (values for currentAveragePrice / LevePrice / LeveProfit columns across ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR sheets, refreshed by the scheduled market-data runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 11377.777
$ws.Range("I43").Value = 2500.5
$ws.Range("J43").Value = 13914.143
$ws.Range("K43").Value = 2500.5
$ws.Range("L43").Value = 13914.143
$ws.Range("M43").Value = -2431.5
$ws.Range("N43").Value = -14052.143

$ws.Range("H47").Value = 11533.5
$ws.Range("I47").Value = 11533.5
$ws.Range("K47").Value = 11533.5
$ws.Range("M47").Value = -10561.5

$ws.Range("H64").Value = 3817.2068
$ws.Range("I64").Value = 3596.3928
$ws.Range("J64").Value = 10000
$ws.Range("K64").Value = 3596.3928
$ws.Range("L64").Value = 10000
$ws.Range("M64").Value = -3348.3928
$ws.Range("N64").Value = -10496

$ws.Range("H67").Value = 3817.2068
$ws.Range("I67").Value = 3596.3928
$ws.Range("J67").Value = 10000
$ws.Range("K67").Value = 3596.3928
$ws.Range("L67").Value = 10000
$ws.Range("M67").Value = -2738.3928
$ws.Range("N67").Value = -11716

$ws.Range("H70").Value = 3000.3333
$ws.Range("J70").Value = 2561.8572
$ws.Range("L70").Value = 7685.571599999999
$ws.Range("N70").Value = -8225.571599999999

$ws.Range("H73").Value = 3000.3333
$ws.Range("J73").Value = 2561.8572
$ws.Range("L73").Value = 7685.571599999999
$ws.Range("N73").Value = -9557.571599999999

$ws.Range("H80").Value = 1779.7354
$ws.Range("I80").Value = 911.1818
$ws.Range("J80").Value = 2195.1304
$ws.Range("K80").Value = 2733.5454
$ws.Range("L80").Value = 6585.3912
$ws.Range("M80").Value = -1735.5454
$ws.Range("N80").Value = -8581.3912

$ws.Range("H83").Value = 1779.7354
$ws.Range("I83").Value = 911.1818
$ws.Range("J83").Value = 2195.1304
$ws.Range("K83").Value = 8200.636199999999
$ws.Range("L83").Value = 19756.1736
$ws.Range("M83").Value = -3208.636199999999
$ws.Range("N83").Value = -29740.1736

$ws.Range("H107").Value = 33956.26
$ws.Range("I107").Value = 554
$ws.Range("J107").Value = 85915.336
$ws.Range("K107").Value = 554
$ws.Range("L107").Value = 85915.336
$ws.Range("M107").Value = 1366
$ws.Range("N107").Value = -89755.336

$ws.Range("H113").Value = 71452.8
$ws.Range("I113").Value = 145470.72
$ws.Range("K113").Value = 145470.72
$ws.Range("M113").Value = -142216.72

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 1500
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H32").Value = 28012.334
$ws.Range("I32").Value = 8238.777
$ws.Range("K32").Value = 8238.777
$ws.Range("M32").Value = -7951.777

$ws.Range("H61").Value = 1893.3182
$ws.Range("I61").Value = 1856
$ws.Range("J61").Value = 2020.2
$ws.Range("K61").Value = 1856
$ws.Range("L61").Value = 2020.2
$ws.Range("M61").Value = -1644
$ws.Range("N61").Value = -2444.2

$ws.Range("H122").Value = 2729.75
$ws.Range("I122").Value = 2734
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 8202
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -5752
$ws.Range("N122").Value = -13000

$ws.Range("H132").Value = 1748.8182
$ws.Range("I132").Value = 1477.3704
$ws.Range("K132").Value = 4432.1112
$ws.Range("M132").Value = -1902.1112

$ws.Range("H136").Value = 1893.3182
$ws.Range("I136").Value = 1856
$ws.Range("J136").Value = 2020.2
$ws.Range("K136").Value = 5568
$ws.Range("L136").Value = 6060.6
$ws.Range("M136").Value = -3018
$ws.Range("N136").Value = -11160.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 855
$ws.Range("I94").Value = 855
$ws.Range("K94").Value = 855
$ws.Range("M94").Value = -404

$ws.Range("H134").Value = 1467.7931
$ws.Range("I134").Value = 1212.3636
$ws.Range("J134").Value = 2270.5715
$ws.Range("K134").Value = 3637.0908
$ws.Range("L134").Value = 6811.7145
$ws.Range("M134").Value = -1102.0908
$ws.Range("N134").Value = -11881.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4313.8
$ws.Range("I16").Value = 5590
$ws.Range("K16").Value = 5590
$ws.Range("M16").Value = -5303

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws.Range("H113").Value = 4313.8
$ws.Range("I113").Value = 5590
$ws.Range("K113").Value = 5590
$ws.Range("M113").Value = -3420

$ws.Range("H122").Value = 104902.7
$ws.Range("I122").Value = 167467.83
$ws.Range("J122").Value = 11055
$ws.Range("K122").Value = 502403.49
$ws.Range("L122").Value = 33165
$ws.Range("M122").Value = -499953.49
$ws.Range("N122").Value = -38065

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

$ws.Range("H132").Value = 4480.8184
$ws.Range("I132").Value = 4473.7646
$ws.Range("J132").Value = 4504.8
$ws.Range("K132").Value = 13421.2938
$ws.Range("L132").Value = 13514.4
$ws.Range("M132").Value = -10891.2938
$ws.Range("N132").Value = -18574.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 134.4
$ws.Range("I7").Value = 159
$ws.Range("K7").Value = 477
$ws.Range("M7").Value = -365

$ws.Range("H37").Value = 500037120
$ws.Range("J37").Value = 500037120
$ws.Range("L37").Value = 1500111360
$ws.Range("N37").Value = -1500111584

$ws.Range("H109").Value = 1359.6
$ws.Range("I109").Value = 1359.6
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 4078.8
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -3038.8
$ws.Range("N109").ClearContents()

$ws.Range("H131").Value = 2044.4546
$ws.Range("I131").Value = 3015
$ws.Range("J131").Value = 1828.7778
$ws.Range("K131").Value = 9045
$ws.Range("L131").Value = 5486.3334
$ws.Range("M131").Value = -4005
$ws.Range("N131").Value = -15566.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3602.2666
$ws.Range("I113").Value = 4131.125
$ws.Range("J113").Value = 2997.8572
$ws.Range("K113").Value = 4131.125
$ws.Range("L113").Value = 2997.8572
$ws.Range("M113").Value = -1961.125
$ws.Range("N113").Value = -7337.8572

$ws.Range("H132").Value = 3429.0527
$ws.Range("I132").Value = 3497.3333
$ws.Range("K132").Value = 10491.9999
$ws.Range("M132").Value = -7961.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1001.8182
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1001.8182
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1001.8182
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1591.8182

$ws.Range("H27").Value = 1001.8182
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1001.8182
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1001.8182
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1215.8182

$ws.Range("H43").Value = 14092590
$ws.Range("J43").Value = 14092590
$ws.Range("L43").Value = 14092590
$ws.Range("N43").Value = -14092976

$ws.Range("H55").Value = 475.2353
$ws.Range("I55").Value = 633.75
$ws.Range("J55").Value = 334.33334
$ws.Range("K55").Value = 633.75
$ws.Range("L55").Value = 334.33334
$ws.Range("M55").Value = -460.75
$ws.Range("N55").Value = -680.33334

$ws.Range("H61").Value = 80761.94
$ws.Range("I61").Value = 92874.91
$ws.Range("K61").Value = 92874.91
$ws.Range("M61").Value = -92672.91

$ws.Range("H69").Value = 44000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 44000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 44000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -45622

$ws.Range("H72").Value = 44000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 44000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 132000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -140112

$ws.Range("H100").Value = 28027.666
$ws.Range("I100").Value = 5693.222
$ws.Range("K100").Value = 5693.222
$ws.Range("M100").Value = -5152.222

$ws.Range("H113").Value = 80761.94
$ws.Range("I113").Value = 92874.91
$ws.Range("K113").Value = 92874.91
$ws.Range("M113").Value = -90704.91

$ws.Range("H125").Value = 80715
$ws.Range("J125").Value = 80715
$ws.Range("L125").Value = 80715
$ws.Range("N125").Value = -90555

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 29000
$ws.Range("J33").Value = 29000
$ws.Range("L33").Value = 29000
$ws.Range("N33").Value = -29500

$ws.Range("H36").Value = 29000
$ws.Range("J36").Value = 29000
$ws.Range("L36").Value = 29000
$ws.Range("N36").Value = -29500

$ws.Range("H62").Value = 8030.6
$ws.Range("I62").Value = 3132.3635
$ws.Range("K62").Value = 3132.3635
$ws.Range("M62").Value = -2508.3635

$ws.Range("H65").Value = 8030.6
$ws.Range("I65").Value = 3132.3635
$ws.Range("K65").Value = 15661.8175
$ws.Range("M65").Value = -12541.8175

$ws.Range("H68").Value = 62758.5
$ws.Range("I68").Value = 25246
$ws.Range("J68").Value = 100271
$ws.Range("K68").Value = 25246
$ws.Range("L68").Value = 100271
$ws.Range("M68").Value = -24435
$ws.Range("N68").Value = -101893

$ws.Range("H71").Value = 62758.5
$ws.Range("I71").Value = 25246
$ws.Range("J71").Value = 100271
$ws.Range("K71").Value = 75738
$ws.Range("L71").Value = 300813
$ws.Range("M71").Value = -71682
$ws.Range("N71").Value = -308925

$ws.Range("H81").Value = 2312.4375
$ws.Range("I81").Value = 2281.6365
$ws.Range("J81").Value = 2380.2
$ws.Range("K81").Value = 4563.273
$ws.Range("L81").Value = 4760.4
$ws.Range("M81").Value = -3502.273
$ws.Range("N81").Value = -6882.4

$ws.Range("H84").Value = 2312.4375
$ws.Range("I84").Value = 2281.6365
$ws.Range("J84").Value = 2380.2
$ws.Range("K84").Value = 22816.365
$ws.Range("L84").Value = 23802
$ws.Range("M84").Value = -17512.365
$ws.Range("N84").Value = -34410

$ws.Range("H113").Value = 402.8889
$ws.Range("I113").Value = 203.375
$ws.Range("K113").Value = 610.125
$ws.Range("M113").Value = 1559.875

$ws.Range("H114").Value = 40000
$ws.Range("J114").Value = 40000
$ws.Range("L114").Value = 40000
$ws.Range("N114").Value = -48678

$ws.Range("H132").Value = 1153.9231
$ws.Range("I132").Value = 777.2286
$ws.Range("K132").Value = 2331.6858
$ws.Range("M132").Value = 198.3141999999998

$ws.Range("H136").Value = 2676.3513
$ws.Range("I136").Value = 2561
$ws.Range("J136").Value = 3035.2222
$ws.Range("K136").Value = 7683
$ws.Range("L136").Value = 9105.6666
$ws.Range("M136").Value = -5133
$ws.Range("N136").Value = -14205.6666

